$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.293.13'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.425.20'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '563.31'
$ws.Range('E5').Value = '  +2.22%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.71'
$ws.Range('E6').Value = '  +3.55%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +1.98%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.423.40'
$ws.Range('E9').Value = '  +1.95%  '
$ws.Range('E10').Value = '  +2.08%  '
$ws.Range('E11').Value = '  -1.56%  '
$ws.Range('E12').Value = '  +0.83%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.354'
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.08'
$ws.Range('E14').Value = '  +2.18%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000178'
$ws.Range('E15').Value = '  +6.03%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.864.37'
$ws.Range('E16').Value = '  +2.11%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.969.46'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.423.58'
$ws.Range('E18').Value = '  +1.99%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.35'
$ws.Range('E19').Value = '  +3.41%  '
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '325.16'
$ws.Range('E21').Value = '  +1.38%  '
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.68'
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('E25').Value = '  -2.41%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.94'
$ws.Range('E26').Value = '  +1.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '590.38'
$ws.Range('E27').Value = '  +15.03%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0₃0947'
$ws.Range('E28').Value = '  +5.72%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.530.55'
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +6.37%  '
$ws.Range('E32').Value = '  +1.42%  '
$ws.Range('E33').Value = '  +0.47%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.89'
$ws.Range('E34').Value = '  +2.65%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.56'
$ws.Range('E35').Value = '  +1.23%  '
$ws.Range('E36').Value = '  +5.24%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.83'
$ws.Range('E38').Value = '  +3.00%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '154.35'
$ws.Range('E39').Value = '  +5.13%  '
$ws.Range('E40').Value = '  +1.63%  '
$ws.Range('E41').Value = '  +1.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.83'
$ws.Range('E42').Value = '  -2.68%  '
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.34'
$ws.Range('E44').Value = '  +8.71%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '150.68'
$ws.Range('E45').Value = '  +1.76%  '
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0540'
$ws.Range('E47').Value = '  +2.70%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '20.47'
$ws.Range('E48').Value = '  +4.38%  '
$ws.Range('E49').Value = '  +2.35%  '
$ws.Range('E50').Value = '  +2.28%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0229'
$ws.Range('E51').Value = '  +1.83%  '
